$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: Name
$ws.Range("B2").Value = "USD"
$ws.Range("C2").Value = "EUR"
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: Yieldcurve
$ws.Range("B3").Value = "USD LIBOR3M OISSTRIPPED 31122019"
$ws.Range("C3").Value = "EURIBOR6M 31122019"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# Row 4: Volatility
$ws.Range("B4").Value = "USD VOL EURUSD 8Y coterminal alpha3perc 31122019"
$ws.Range("C4").Value = "EUR VOL EURUSD 8Y coterminal alpha3perc 31122019"
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Row 5: MeanReversion - keep B5, C5 as 0.03 ; clear D5,E5 but keep formatting (percentage style)
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()

# Update selection to D11
$ws.Range("D11").Select()
